# se modifica data para empezar regresion R34 en Pre Prod
#
# Updates the regression-seed rows (2 and 3) on Hoja1 from the "JunioUno"
# batch to the "JulioUno" batch, bumps the associated date/ID counters, and
# scrolls the sheet so column G is visible at the left edge (topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2: Answ/Ans RegrJunioUno -> Answ/Ans RegrJulioUno, date + id bump
$ws.Range("F2").Value = "AnswRegrJulioUno"
$ws.Range("G2").Value = "AnsRegrJulioUno"
$ws.Range("H2").Value = 21300129
$ws.Range("O2").Value = 128

# Row 3: date + id bump (text values for F3/G3 are unchanged)
$ws.Range("H3").Value = 21300130
$ws.Range("O3").Value = 129

# Scroll the sheet view so that column G is the left-most visible column
# (mirrors <sheetView topLeftCell="G1" .../> in the saved workbook).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1

$wb.Save()
